# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated cryptos list on Tue Jan  9 14:42:47 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, [string]$val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "46.686.29"
Set-TextValue $ws "E2" "  +3.74%  "
Set-TextValue $ws "D3" "2.258.65"
Set-TextValue $ws "E3" "  -0.38%  "
Set-TextValue $ws "E4" "  -0.02%  "
Set-TextValue $ws "D5" "301.11"
Set-TextValue $ws "E5" "  -0.41%  "
Set-TextValue $ws "D6" "100.38"
Set-TextValue $ws "E6" "  +5.51%  "
Set-TextValue $ws "E7" "  -0.50%  "
Set-TextValue $ws "E8" "  +0.04%  "
Set-TextValue $ws "D9" "0.512"
Set-TextValue $ws "E9" "  +0.53%  "
Set-TextValue $ws "D10" "35.54"
Set-TextValue $ws "E10" "  +3.73%  "
Set-TextValue $ws "D11" "0.0781"
Set-TextValue $ws "E11" "  -1.10%  "
Set-TextValue $ws "D12" "7.16"
Set-TextValue $ws "E12" "  -0.53%  "
Set-TextValue $ws "D13" "0.102"
Set-TextValue $ws "E13" "  -0.79%  "
Set-TextValue $ws "D14" "2.610.87"
Set-TextValue $ws "E14" "  -0.18%  "
Set-TextValue $ws "D15" "2.265.75"
Set-TextValue $ws "E15" "  -0.15%  "
Set-TextValue $ws "D16" "13.54"
Set-TextValue $ws "E16" "  -0.72%  "
Set-TextValue $ws "D17" "46.667.85"
Set-TextValue $ws "E17" "  +3.90%  "
Set-TextValue $ws "D18" "0.794"
Set-TextValue $ws "E18" "  -0.60%  "
Set-TextValue $ws "D19" "12.93"
Set-TextValue $ws "E19" "  -0.02%  "
Set-TextValue $ws "D20" "0.0₃0927"
Set-TextValue $ws "E20" "  +0.46%  "
Set-TextValue $ws "D21" "5.89"
Set-TextValue $ws "E21" "  -3.23%  "
Set-TextValue $ws "D22" "65.31"
Set-TextValue $ws "E22" "  -0.40%  "
Set-TextValue $ws "D23" "249.27"
Set-TextValue $ws "E23" "  +4.57%  "
Set-TextValue $ws "D24" "2.83"
Set-TextValue $ws "E24" "  -2.22%  "
Set-TextValue $ws "E25" "  +0.07%  "
Set-TextValue $ws "D26" "1.87"
Set-TextValue $ws "E26" "  -0.99%  "
Set-TextValue $ws "D27" "42.96"
Set-TextValue $ws "E27" "  +3.97%  "
Set-TextValue $ws "D28" "2.24"
Set-TextValue $ws "E28" "  -2.54%  "
Set-TextValue $ws "D29" "9.72"
Set-TextValue $ws "E29" "  +1.35%  "
Set-TextValue $ws "D30" "19.83"
Set-TextValue $ws "E30" "  +1.36%  "
Set-TextValue $ws "E31" "  +8.58%  "
Set-TextValue $ws "D32" "146.69"
Set-TextValue $ws "E32" "  -4.27%  "
Set-TextValue $ws "D33" "5.43"
Set-TextValue $ws "E33" "  -2.47%  "
Set-TextValue $ws "B34" "Hedera"
Set-TextValue $ws "C34" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D34" "0.0769"
Set-TextValue $ws "E34" "  -2.26%  "
Set-TextValue $ws "B35" "LidoDAOToken"
Set-TextValue $ws "C35" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws "D35" "3.18"
Set-TextValue $ws "E35" "  +7.62%  "
Set-TextValue $ws "E36" "  +9.56%  "
Set-TextValue $ws "D37" "0.115"
Set-TextValue $ws "E37" "  -1.20%  "
Set-TextValue $ws "D38" "16.29"
Set-TextValue $ws "E38" "  +19.66%  "
Set-TextValue $ws "D39" "1.72"
Set-TextValue $ws "E39" "  -2.81%  "
Set-TextValue $ws "D40" "3.86"
Set-TextValue $ws "E40" "  -5.11%  "
Set-TextValue $ws "D41" "0.0297"
Set-TextValue $ws "E41" "  -5.10%  "
Set-TextValue $ws "D42" "3.21"
Set-TextValue $ws "E42" "  -1.83%  "
Set-TextValue $ws "E43" "  -0.16%  "
Set-TextValue $ws "D44" "1.98"
Set-TextValue $ws "E44" "  +1.64%  "
Set-TextValue $ws "D45" "1.809.62"
Set-TextValue $ws "E45" "  +3.28%  "
Set-TextValue $ws "D46" "90.44"
Set-TextValue $ws "E46" "  +19.08%  "
Set-TextValue $ws "D47" "72.52"
Set-TextValue $ws "E47" "  +2.46%  "
Set-TextValue $ws "D48" "0.188"
Set-TextValue $ws "E48" "  -4.75%  "
Set-TextValue $ws "D49" "4.81"
Set-TextValue $ws "E49" "  +2.69%  "
Set-TextValue $ws "D50" "93.68"
Set-TextValue $ws "E50" "  -2.82%  "
Set-TextValue $ws "D51" "2.487.63"
Set-TextValue $ws "E51" "  -0.15%  "
